# feat: add 2022-Q3 data
#
# 1. "总计" (total) sheet: insert a new row right below the header for the
#    2022-Q3 summary figures, pushing the existing 2021-Q4 / 2021-Q1 rows
#    down by one.
# 2. Insert a brand-new worksheet named "2022-Q3" (placed right after the
#    "总计" sheet, i.e. as the new 2nd sheet) holding the per-fund detail
#    table, mirroring the layout of the existing quarterly sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: "总计" sheet - insert the 2022-Q3 row
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Push rows 2-3 down to make room for the new 2022-Q3 row.
$total.Range("A2").EntireRow.Insert()

# The Insert() above copies the formatting of the row above (the header),
# which is not what we want for a plain data row - clear it first …
$total.Range("B2:D2").ClearFormats()
# … then copy the real target formatting (matching the other data rows)
# from row 3 into row 2's A cell (only column A carries an explicit style).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.05

# The pandas-style row index in column A is renumbered sequentially after
# the insert (0, 1, 2, ...), so refresh the two pushed-down rows as well.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

# ---------------------------------------------------------------------
# Step 2: brand-new "2022-Q3" worksheet with the fund detail table
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# Header row
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"
$q3.Range("B1:H1").Font.Bold = $true

# Make sure the text-looking numeric columns are stored as plain text
# (matching the source data, e.g. fund codes with leading zeros like
# "014781", and figures such as "1.77").
$q3.Range("B2:G3").NumberFormat = "@"

# Row 2 - 建信兴衡优选一年持有期混合A
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "014781"
$q3.Range("C2").Value = "建信兴衡优选一年持有期混合A"
$q3.Range("D2").Value = "1.77"
$q3.Range("E2").Value = "46.74"
$q3.Range("F2").Value = "1.83"
$q3.Range("G2").Value = "0.0324"
$q3.Range("H2").Value = 10

# Row 3 - 建信兴衡优选一年持有期混合C
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "014782"
$q3.Range("C3").Value = "建信兴衡优选一年持有期混合C"
$q3.Range("D3").Value = "0.70"
$q3.Range("E3").Value = "46.74"
$q3.Range("F3").Value = "1.83"
$q3.Range("G3").Value = "0.0128"
$q3.Range("H3").Value = 10

# Restore the original active sheet ("2021-Q1" - the last sheet) as the
# selected / active tab, since adding a new sheet makes it active.
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
